# AESTHEMOS_Items_Naming workbook update
# - Adds a third "Raw Data" column (with the original "AeNN" item codes) to the
#   query table / worksheet on the "Sheet2" tab.
# - Corrects five "Item Text" values that were re-worded.
# - Updates the selected cell to match the author's last position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet2")

# --- 1. Grow the table by one column and give it a header -------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null
$ws.Range("C1").Value = "Raw Data"

# --- 2. Fill in the new "Raw Data" column (rows 2-43) ------------------------
$rawData = @(
    "Ae25", "Ae17", "Ae31", "Ae40", "Ae14", "Ae36", "Ae33", "Ae19", "Ae24",
    "Ae37", "Ae8",  "Ae18", "Ae41", "Ae16", "Ae34", "Ae7",  "Ae6",  "Ae1",
    "Ae35", "Ae12", "Ae22", "Ae42", "Ae21", "Ae13", "Ae2",  "Ae10", "Ae5",
    "Ae38", "Ae3",  "Ae39", "Ae28", "Ae26", "Ae4",  "Ae20", "Ae23", "Ae15",
    "Ae11", "Ae29", "Ae27", "Ae30", "Ae32", "Ae9"
)

for ($i = 0; $i -lt $rawData.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $rawData[$i]
}

# --- 3. Re-word five "Item Text" cells in column B ---------------------------
$ws.Range("B10").Value = "Made me feel confused"
$ws.Range("B12").Value = "I felt Something Wonderful"
$ws.Range("B13").Value = "Was Enchanting"
$ws.Range("B16").Value = "Was Impressive"
$ws.Range("B21").Value = "Was ugly"

# --- 4. Cosmetic touch-ups ----------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 11.43
$ws.Range("E41").Select()
